$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row for the section header "grandes regiões e unidades da federação"
# (row 6) is removed entirely; all rows below it shift up by one, and the
# now-unused shared string is dropped automatically when the workbook is
# saved.
$ws.Rows.Item(6).Delete()
